# Ajout du dossier de conception dans le rapport de projet
# Append two new rows (24 and 25) to the "Journal" sheet's Tableau1,
# mirroring the formatting of the last existing row (23), then extend the
# table range to cover them and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (number formats / alignment) of the last data row
# onto the two new rows before writing values into them.
$ws.Range("A23:E23").Copy()
$ws.Range("A24:E25").PasteSpecial(-4122)

# Row 24: Documentation entry about the "dossier de conception".
$ws.Range("A24").Value = 44981
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = 0.041666666666666664
$ws.Range("D24").Value = "Documentation"
$ws.Range("E24").Value = "Rédaction du dossier de conception dans le rapport de projet"

# Row 25: Documentation entry about the "résumé de projet" (no hours logged yet).
$ws.Range("A25").Value = 44981
$ws.Range("B25").Value = 3
$ws.Range("C25").ClearContents()
$ws.Range("D25").Value = "Documentation"
$ws.Range("E25").Value = "Rédaction du résumé de projet "

# Grow Tableau1 so the two new rows become part of the table.
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("A1:E25"))

# Match the author's final selection.
$ws.Range("E26").Select()
